# Updates cryptos list data (price & volume columns, plus the Uniswap/Polygon
# row swap) to match the latest scrape, per commit:
# "Updated cryptos list on Sat May 18 15:00:19 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.809.94"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.104.80"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.34"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.49"
$ws.Range("E6").Value = "  +2.82%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.100.64"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.43"
$ws.Range("E10").Value = "  -3.71%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.479"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.16"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.618.89"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.797.94"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.10"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.103.12"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.37"
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "476.41"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.00"
$ws.Range("E22").Value = "  +7.30%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.712"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.39"
$ws.Range("E24").Value = "  +4.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.78"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.88"
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.42"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.89"
$ws.Range("E30").Value = "  -2.15%  "
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.61"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0930"
$ws.Range("E34").Value = "  -9.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.973"
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.42"
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.08"
$ws.Range("E39").Value = "  +1.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.96"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.309"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.58"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.796.11"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0354"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "377.54"
$ws.Range("E46").Value = "  -2.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.54"
$ws.Range("E47").Value = "  -12.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.89"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.83"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("E51").Value = "  -1.07%  "
